$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a "Price" cell (column D). These are stored as text in the
# sheet. Values with a single decimal point read back as numbers unless
# forced to text with a leading apostrophe (values with multiple '.'
# thousands separators are never numeric, so they don't need the prefix),
# mirroring how the source data was authored.
function Set-Price($row, $value) {
    if ($value -match '^[0-9]+(\.[0-9]+)?$') {
        $ws.Range("D$row").Value = "'" + $value
    } else {
        $ws.Range("D$row").Value = $value
    }
}
function Set-Volume($row, $value) {
    $ws.Range("E$row").Value = "  $value  "
}

# Row 2 - Bitcoin
Set-Price 2 "51.894.56"
Set-Volume 2 "+0.13%"

# Row 3 - Ethereum
Set-Price 3 "2.934.88"
Set-Volume 3 "+3.44%"

# Row 4 - TetherUSD
Set-Volume 4 "-0.03%"

# Row 5 - BNB
Set-Price 5 "352.44"
Set-Volume 5 "+0.17%"

# Row 6 - Solana
Set-Price 6 "112.05"
Set-Volume 6 "-0.97%"

# Row 7 - XRP
Set-Price 7 "0.561"
Set-Volume 7 "+0.22%"

# Row 8 - USDC
Set-Volume 8 "+0.00%"

# Row 9 - Cardano
Set-Price 9 "0.620"
Set-Volume 9 "-0.55%"

# Row 10 - Avalanche
Set-Price 10 "39.42"
Set-Volume 10 "-2.03%"

# Row 11 - Dogecoin
Set-Price 11 "0.0887"
Set-Volume 11 "+4.16%"

# Row 12 - TRON
Set-Volume 12 "+1.18%"

# Row 13 - Chainlink
Set-Price 13 "20.11"
Set-Volume 13 "+0.24%"

# Row 14 - Polkadot
Set-Price 14 "7.79"
Set-Volume 14 "-0.43%"

# Row 15 - WrappedliquidstakedEther2.0
Set-Price 15 "3.391.02"
Set-Volume 15 "+3.27%"

# Row 16 - WrappedEther
Set-Price 16 "2.926.02"
Set-Volume 16 "+2.69%"

# Row 17 - Polygon
Set-Volume 17 "+0.65%"

# Row 18 - WrappedBTC
Set-Price 18 "51.976.38"
Set-Volume 18 "+0.14%"

# Row 19 - ImmutableX
Set-Price 19 "3.32"
Set-Volume 19 "-4.31%"

# Row 20 - Uniswap
Set-Price 20 "7.62"
Set-Volume 20 "-0.04%"

# Row 21 - InternetComputer(DFINITY)
Set-Price 21 "14.24"
Set-Volume 21 "+6.46%"

# Row 22 - ShibaInu
Set-Volume 22 "+0.90%"

# Row 23 - Litecoin
Set-Price 23 "71.25"
Set-Volume 23 "+0.96%"

# Row 24 - BitcoinCash
Set-Price 24 "268.49"
Set-Volume 24 "-0.30%"

# Row 25 - PancakeSwap
Set-Price 25 "2.78"
Set-Volume 25 "+0.09%"

# Row 26 - Kaspa
Set-Price 26 "0.181"
Set-Volume 26 "+11.67%"

# Row 27 - EthereumClassic
Set-Price 27 "26.98"
Set-Volume 27 "+2.43%"

# Row 28 - Dai
Set-Price 28 "0.999"
Set-Volume 28 "-0.12%"

# Row 29 - Filecoin
Set-Price 29 "7.23"
Set-Volume 29 "+14.22%"

# Row 30 - Hedera
Set-Price 30 "0.104"
Set-Volume 30 "+15.59%"

# Row 31 - Cosmos
Set-Price 31 "10.57"
Set-Volume 31 "-0.36%"

# Row 32 - now RenderToken (was InjectiveProtocol)
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-Price 32 "6.25"
Set-Volume 32 "+10.31%"

# Row 33 - Toncoin
Set-Volume 33 "+0.06%"

# Row 34 - now InjectiveProtocol (was RenderToken)
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-Price 34 "37.17"
Set-Volume 34 "-4.51%"

# Row 35 - OKB
Set-Price 35 "52.92"
Set-Volume 35 "+0.17%"

# Row 36 - VeChain
Set-Price 36 "0.0454"
Set-Volume 36 "-0.56%"

# Row 37 - FirstDigitalUSD
Set-Price 37 "0.998"
Set-Volume 37 "-0.18%"

# Row 38 - LidoDAOToken
Set-Volume 38 "+3.45%"

# Row 39 - Celestia
Set-Volume 39 "-1.67%"

# Row 40 - ARBITRUM
Set-Volume 40 "+1.40%"

# Row 41 - Stacks
Set-Volume 41 "+6.32%"

# Row 42 - Stellar
Set-Price 42 "0.118"
Set-Volume 42 "+1.19%"

# Row 43 - EnergySwap
Set-Price 43 "23.31"
Set-Volume 43 "+5.02%"

# Row 44 - WEMIXToken
Set-Volume 44 "-1.14%"

# Row 45 - ApeXProtocol
Set-Volume 45 "+2.13%"

# Row 46 - NEARProtocol
Set-Volume 46 "-0.06%"

# Row 47 - Maker
Set-Price 47 "2.174.19"
Set-Volume 47 "-0.51%"

# Row 48 - Monero
Set-Price 48 "111.72"
Set-Volume 48 "-8.52%"

# Row 49 - TheGraph
Set-Volume 49 "+2.66%"

# Row 50 - BEAM
Set-Price 50 "0.0351"
Set-Volume 50 "+11.19%"

# Row 51 - SEI
Set-Price 51 "0.946"
Set-Volume 51 "-2.54%"
